# A new weekly Papaya price-report record (Vega Modelo de Temuco) needs to be
# inserted at the top of the data table (row 49, right after the header block
# that precedes this product's rows). All existing records from row 49 down
# get pushed one row lower (49->50, 50->51, ... 104->105), and the new record
# is written into the freed-up row 49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 49:104 down to 50:105, leaving a blank row 49 behind.
$ws.Rows("49:49").Insert()

# Populate the newly inserted row 49 with the new observation.
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 45096
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100108
$ws.Range("H49").Value = "Tropicales y subtropicales"
$ws.Range("I49").Value = 100108004
$ws.Range("J49").Value = "Papaya"
$ws.Range("K49").Value = "Cultivar IV Región"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 65
$ws.Range("N49").Value = 27000
$ws.Range("O49").Value = 27000
$ws.Range("P49").Value = 27000
$ws.Range("Q49").Value = "$/bandeja 10 kilos"
$ws.Range("R49").Value = "Provincia del Elquí"
$ws.Range("S49").Value = 2700
$ws.Range("T49").Value = 10
